$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "comorbidities_smoking_hx: Y"; New = "comorbidities_smoking_hx" },
    @{ Old = "comorbidities_DM: Y"; New = "comorbidities_DM" },
    @{ Old = "comorbidities_prior_cardiac_surg: Y"; New = "comorbidities_prior_cardiac_surg" },
    @{ Old = "preop_temp_MCS: Y"; New = "preop_temp_MCS" },
    @{ Old = "preop_IABP: Y"; New = "preop_IABP" },
    @{ Old = "preop_imeplla5.5: Y"; New = "preop_imeplla5.5" },
    @{ Old = "preop_VA_ECMO: Y"; New = "preop_VA_ECMO" },
    @{ Old = "preop_LVAD: Y"; New = "preop_LVAD" },
    @{ Old = "rx_preop_inotrope: Y"; New = "rx_preop_inotrope" },
    @{ Old = "rx_preop_amiodarone: Y"; New = "rx_preop_amiodarone" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
